$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.630.28"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.138.82"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.011"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.64%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "352.52"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +5.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.010"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5272"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4553"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.96"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09120"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.182"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.83"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.87%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.129.01"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.856"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.122"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "102.15"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +5.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001172"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.012"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06723"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.50"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.45%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.343"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.701.26"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.79"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.384"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.360.89"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.54"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.646"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "164.82"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "137.02"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.218"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.76%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.667"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.348"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.020"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.203"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +7.42%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.48"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02651"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06920"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2328"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6973"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.274"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.82"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +5.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.349"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6457"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.46%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +5.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.758"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.256"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "83.09"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07296"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.26%  "
